# Updates stock counts/values for several items (apparent stock correction /
# reconciliation) and refreshes the dependent Sub Total / Grand Total cells
# on the single worksheet of the Companywise Stock Report workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 138
$ws.Range("G6").Value = 4123.44
$ws.Range("B10").Value = 35774.57
$ws.Range("F93").Value = 70
$ws.Range("G93").Value = 5528.6
$ws.Range("F94").Value = 146
$ws.Range("G94").Value = 20491.1
$ws.Range("F103").Value = 10
$ws.Range("G103").Value = 508
$ws.Range("F109").Value = 101
$ws.Range("G109").Value = 12672.47
$ws.Range("B114").Value = 230615.02
$ws.Range("F214").Value = 24
$ws.Range("G214").Value = 6888.24
$ws.Range("B222").Value = 41923.69
$ws.Range("F236").Value = 83
$ws.Range("G236").Value = 3566.51
$ws.Range("B274").Value = 63743.38
$ws.Range("B277").Value = 63565
$ws.Range("E277").Value = 109.19
$ws.Range("F277").Value = 60
$ws.Range("G277").Value = 6162.6
$ws.Range("B278").Value = 61610
$ws.Range("E278").Value = 122.71
$ws.Range("F278").Value = -58
$ws.Range("G278").Value = -5957.18
$ws.Range("B294").Value = 63531
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 79
$ws.Range("G294").Value = 11334.92
$ws.Range("B295").Value = 57802
$ws.Range("E295").Value = 162.71
$ws.Range("F295").Value = -79
$ws.Range("G295").Value = -11334.92
$ws.Range("B296").Value = 63571
$ws.Range("F296").Value = 0
$ws.Range("G296").Value = 0
$ws.Range("F301").Value = 51
$ws.Range("G301").Value = 2617.32
$ws.Range("F304").Value = 14
$ws.Range("G304").Value = 4239.06
$ws.Range("F310").Value = 18
$ws.Range("G310").Value = 938.34
$ws.Range("B339").Value = 251591.74
$ws.Range("B356").Value = 63681
$ws.Range("E356").Value = 23.84
$ws.Range("F356").Value = 0
$ws.Range("G356").Value = 0
$ws.Range("B357").Value = 31930
$ws.Range("E357").Value = 26.8
$ws.Range("F357").Value = -62
$ws.Range("G357").Value = -1390.04
$ws.Range("F368").Value = 152
$ws.Range("G368").Value = 25709.28
$ws.Range("F379").Value = 69
$ws.Range("G379").Value = 6416.31
$ws.Range("F393").Value = 94
$ws.Range("G393").Value = 19823.66
$ws.Range("B395").Value = 224054.87
$ws.Range("F442").Value = 25
$ws.Range("G442").Value = 1273.75
$ws.Range("B448").Value = 35305.46
$ws.Range("B465").Value = 53757
$ws.Range("E465").Value = 16.08
$ws.Range("F465").Value = -159
$ws.Range("G465").Value = -2138.55
$ws.Range("B466").Value = 65069
$ws.Range("E466").Value = 14.3
$ws.Range("F466").Value = 2
$ws.Range("G466").Value = 26.9
$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 119
$ws.Range("G479").Value = 1930.18
$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68
$ws.Range("F484").Value = 378
$ws.Range("G484").Value = 2453.22
$ws.Range("B485").Value = 64925
$ws.Range("E485").Value = 13.97
$ws.Range("F485").Value = 129
$ws.Range("G485").Value = 1696.35
$ws.Range("B486").Value = 45709
$ws.Range("E486").Value = 15.69
$ws.Range("F486").Value = -300
$ws.Range("G486").Value = -3945
$ws.Range("B490").Value = 53595
$ws.Range("E490").Value = 17.61
$ws.Range("F490").Value = -335
$ws.Range("G490").Value = -4934.55
$ws.Range("B491").Value = 65067
$ws.Range("E491").Value = 15.65
$ws.Range("F491").Value = 172
$ws.Range("G491").Value = 2533.56
$ws.Range("B492").Value = -14067.84
$ws.Range("F541").Value = 15
$ws.Range("G541").Value = 1654.05
$ws.Range("B542").Value = 6730.87
$ws.Range("F548").Value = 193
$ws.Range("G548").Value = 1312.4
$ws.Range("B556").Value = 5009.57
$ws.Range("B568").Value = 64810
$ws.Range("E568").Value = 291.22
$ws.Range("F568").Value = 5
$ws.Range("G568").Value = 1369.6
$ws.Range("B569").Value = 53319
$ws.Range("E569").Value = 310.64
$ws.Range("F569").Value = -6
$ws.Range("G569").Value = -1643.52
$ws.Range("F580").Value = 68
$ws.Range("G580").Value = 2322.2
$ws.Range("F584").Value = 38
$ws.Range("G584").Value = 689.3200000000001
$ws.Range("B586").Value = 18388.74
$ws.Range("B600").Value = 60022
$ws.Range("E600").Value = 37.22
$ws.Range("F600").Value = -113
$ws.Range("G600").Value = -3709.79
$ws.Range("B601").Value = 64830
$ws.Range("E601").Value = 34.9
$ws.Range("F601").Value = 111
$ws.Range("G601").Value = 3644.13
$ws.Range("F622").Value = 75
$ws.Range("G622").Value = 9791.25
$ws.Range("B625").Value = 17571.76
$ws.Range("F654").Value = 296
$ws.Range("G654").Value = 23792.48
$ws.Range("B655").Value = 31691.92
$ws.Range("F705").Value = 76
$ws.Range("G705").Value = 10877.88
$ws.Range("B709").Value = 61428
$ws.Range("D709").Value = 69.16
$ws.Range("E709").Value = 73.52
$ws.Range("F709").Value = 1
$ws.Range("G709").Value = 69.16
$ws.Range("B710").Value = 63150
$ws.Range("D710").Value = 75.68000000000001
$ws.Range("E710").Value = 80.45
$ws.Range("F710").Value = 20
$ws.Range("G710").Value = 1513.6
$ws.Range("F714").Value = 59
$ws.Range("G714").Value = 4104.04
$ws.Range("F717").Value = 158
$ws.Range("G717").Value = 21331.58
$ws.Range("B720").Value = 58183.69
$ws.Range("B736").Value = 65362
$ws.Range("F736").Value = 37
$ws.Range("G736").Value = 1512.19
$ws.Range("B737").Value = 65079
$ws.Range("F737").Value = 21
$ws.Range("G737").Value = 858.27
$ws.Range("B747").Value = 51637.7
$ws.Range("F772").Value = 2626
$ws.Range("G772").Value = 428326.86
$ws.Range("B779").Value = 651688.84
$ws.Range("F784").Value = 20
$ws.Range("G784").Value = 819.8
$ws.Range("F790").Value = 4
$ws.Range("G790").Value = 642.04
$ws.Range("B796").Value = 61663.16
$ws.Range("B797").Value = 2434310.3
$ws.Range("B798").Value = 2434310.3
